$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '30.428.23'
$ws.Cells.Item(2, 5).Value = '  -0.88%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.917.80'
$ws.Cells.Item(3, 5).Value = '  +2.08%  '

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9994'
$cell.ClearFormats()
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '241.79'
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +1.67%  '

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9993'
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -0.03%  '

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4703'
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -0.99%  '

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.2849'
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  +0.72%  '

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06834'
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +4.88%  '

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '108.41'
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +13.15%  '

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.29'
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -1.84%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.894.04'
$ws.Cells.Item(12, 5).Value = '  +0.64%  '

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07647'
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +0.84%  '

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.212'
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  +3.04%  '

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.6569'
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  +1.15%  '

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '289.37'
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  -4.50%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '30.433.75'
$ws.Cells.Item(17, 5).Value = '  -0.80%  '

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.000007641'
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = '  +1.60%  '

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9990'
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  -0.16%  '

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.93'
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -0.84%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '2.145.59'
$ws.Cells.Item(21, 5).Value = '  +0.91%  '

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9999'
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +0.07%  '

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.225'
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +1.75%  '

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.202'
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +0.85%  '

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '21.58'
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +9.76%  '

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '167.91'
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -0.92%  '

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.281'
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +0.33%  '

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.051'
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +5.13%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1072'
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +1.35%  '

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.371'
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +1.50%  '

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.156'
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  +0.05%  '

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.960'
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +0.24%  '

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05051'
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +0.87%  '

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7417'
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  +3.15%  '

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.154'
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -1.34%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'VeChain'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.02081'
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  +8.74%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'HuobiToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.750'
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  +1.55%  '

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.694'
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -0.04%  '

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.056'
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +0.44%  '

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.8778'
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -2.22%  '

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '108.74'
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +1.54%  '

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.867'
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +5.05%  '

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9990'
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  -0.10%  '

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4226'
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +0.93%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '67.55'
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +3.47%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'BitcoinSV'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '50.59'
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = '  +18.39%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.170'
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -2.10%  '

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.201'
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +2.42%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Algorand'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1210'
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -0.26%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Elrond'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '34.78'
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  +0.89%  '

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3895'
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +2.69%  '
